# Updates cryptos list price/volume figures (cap3k GitHub Actions refresh).
# Price (D) and Volume(1h) (E) columns are stored as plain text in this
# sheet (e.g. "61.454.00", "0.999", "  +0.61%  "), so a bare Range.Value
# assignment would let Excel auto-coerce numeric-looking strings into real
# numbers. Force-format the cell as Text first, assign, then restore the
# original style so no visible/persisted formatting changes leak in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '61.454.00'
Set-TextValue $ws.Range('E2') '  +0.61%  '
Set-TextValue $ws.Range('D3') '2.933.33'
Set-TextValue $ws.Range('E3') '  +0.18%  '
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '598.77'
Set-TextValue $ws.Range('E5') '  +1.37%  '
Set-TextValue $ws.Range('D6') '145.27'
Set-TextValue $ws.Range('E6') '  -0.16%  '
Set-TextValue $ws.Range('E7') '  -0.04%  '
Set-TextValue $ws.Range('E8') '  -0.72%  '
Set-TextValue $ws.Range('D9') '6.98'
Set-TextValue $ws.Range('E9') '  +1.76%  '
Set-TextValue $ws.Range('D11') '0.440'
Set-TextValue $ws.Range('E11') '  -0.77%  '
Set-TextValue $ws.Range('E12') '  -0.75%  '
Set-TextValue $ws.Range('D13') '33.68'
Set-TextValue $ws.Range('E13') '  -0.13%  '
Set-TextValue $ws.Range('E14') '  +0.55%  '
Set-TextValue $ws.Range('D15') '3.419.30'
Set-TextValue $ws.Range('E15') '  +0.28%  '
Set-TextValue $ws.Range('D16') '61.466.96'
Set-TextValue $ws.Range('E16') '  +0.73%  '
Set-TextValue $ws.Range('D17') '6.72'
Set-TextValue $ws.Range('E17') '  -0.20%  '
Set-TextValue $ws.Range('D18') '2.931.19'
Set-TextValue $ws.Range('E18') '  +0.21%  '
Set-TextValue $ws.Range('D19') '432.96'
Set-TextValue $ws.Range('E19') '  +0.32%  '
Set-TextValue $ws.Range('D20') '13.50'
Set-TextValue $ws.Range('E20') '  +0.07%  '
Set-TextValue $ws.Range('D21') '0.679'
Set-TextValue $ws.Range('E21') '  -0.83%  '
Set-TextValue $ws.Range('D22') '7.12'
Set-TextValue $ws.Range('E22') '  +0.15%  '
Set-TextValue $ws.Range('D23') '81.87'
Set-TextValue $ws.Range('E23') '  +0.87%  '
Set-TextValue $ws.Range('D24') '10.92'
Set-TextValue $ws.Range('E24') '  -1.38%  '
Set-TextValue $ws.Range('D25') '2.20'
Set-TextValue $ws.Range('E25') '  -1.08%  '
Set-TextValue $ws.Range('D26') '11.79'
Set-TextValue $ws.Range('E26') '  -2.07%  '
Set-TextValue $ws.Range('E27') '  -0.08%  '
Set-TextValue $ws.Range('D28') '2.24'
Set-TextValue $ws.Range('E28') '  -3.55%  '
Set-TextValue $ws.Range('D29') '2.62'
Set-TextValue $ws.Range('E29') '  -0.17%  '
Set-TextValue $ws.Range('D30') '6.93'
Set-TextValue $ws.Range('D31') '26.73'
Set-TextValue $ws.Range('E32') '  +1.12%  '
Set-TextValue $ws.Range('D33') '0.999'
Set-TextValue $ws.Range('E33') '  -0.03%  '
Set-TextValue $ws.Range('D34') '0.0₃0881'
Set-TextValue $ws.Range('E34') '  +1.92%  '
Set-TextValue $ws.Range('E35') '  -0.07%  '
Set-TextValue $ws.Range('D36') '5.65'
Set-TextValue $ws.Range('E36') '  +0.07%  '
Set-TextValue $ws.Range('D37') '3.01'
Set-TextValue $ws.Range('E37') '  -1.91%  '
Set-TextValue $ws.Range('D38') '2.01'
Set-TextValue $ws.Range('E38') '  +0.10%  '
Set-TextValue $ws.Range('E39') '  -0.81%  '
Set-TextValue $ws.Range('D40') '8.62'
Set-TextValue $ws.Range('E40') '  -0.01%  '
Set-TextValue $ws.Range('D41') '42.55'
Set-TextValue $ws.Range('E41') '  +8.00%  '
Set-TextValue $ws.Range('D42') '0.283'
Set-TextValue $ws.Range('E42') '  -1.29%  '
Set-TextValue $ws.Range('E43') '  -0.13%  '
Set-TextValue $ws.Range('D44') '2.704.55'
Set-TextValue $ws.Range('E44') '  -0.24%  '
Set-TextValue $ws.Range('D45') '134.57'
Set-TextValue $ws.Range('E45') '  +2.02%  '
Set-TextValue $ws.Range('D46') '364.84'
Set-TextValue $ws.Range('E46') '  -3.16%  '
Set-TextValue $ws.Range('D48') '23.78'
Set-TextValue $ws.Range('E48') '  -1.87%  '
Set-TextValue $ws.Range('E49') '  -1.44%  '
Set-TextValue $ws.Range('E50') '  -1.75%  '
Set-TextValue $ws.Range('E51') '  -0.82%  '
